$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 11 (Marking): Right 5 -> 4, Wrong -1 -> -2
$ws.Range("B11").Value = 4
$ws.Range("C11").Value = -2

# Row 12 (Total): Right 95 -> 76, Wrong -4 -> -8, Max text "95 / 140" -> "68 / 112"
$ws.Range("B12").Value = 76
$ws.Range("C12").Value = -8
$ws.Range("E12").Value = "68 / 112"
